$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormatLocal = "@"
$ws.Range("D2").Value = "64.930.65"
$ws.Range("D2").NumberFormatLocal = "General"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").NumberFormatLocal = "@"
$ws.Range("D3").Value = "3.394.83"
$ws.Range("D3").NumberFormatLocal = "General"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormatLocal = "@"
$ws.Range("D5").Value = "561.08"
$ws.Range("D5").NumberFormatLocal = "General"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormatLocal = "@"
$ws.Range("D6").Value = "175.18"
$ws.Range("D6").NumberFormatLocal = "General"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("D7").NumberFormatLocal = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").NumberFormatLocal = "General"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").NumberFormatLocal = "@"
$ws.Range("D8").Value = "3.385.49"
$ws.Range("D8").NumberFormatLocal = "General"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormatLocal = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").NumberFormatLocal = "General"
$ws.Range("E10").Value = "  +9.50%  "
$ws.Range("D11").NumberFormatLocal = "@"
$ws.Range("D11").Value = "0.633"
$ws.Range("D11").NumberFormatLocal = "General"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").NumberFormatLocal = "@"
$ws.Range("D12").Value = "54.66"
$ws.Range("D12").NumberFormatLocal = "General"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").NumberFormatLocal = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").NumberFormatLocal = "General"
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").NumberFormatLocal = "@"
$ws.Range("D15").Value = "3.936.42"
$ws.Range("D15").NumberFormatLocal = "General"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").NumberFormatLocal = "@"
$ws.Range("D16").Value = "18.29"
$ws.Range("D16").NumberFormatLocal = "General"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").NumberFormatLocal = "@"
$ws.Range("D17").Value = "3.386.26"
$ws.Range("D17").NumberFormatLocal = "General"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormatLocal = "@"
$ws.Range("D19").Value = "11.92"
$ws.Range("D19").NumberFormatLocal = "General"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormatLocal = "@"
$ws.Range("D20").Value = "64.804.31"
$ws.Range("D20").NumberFormatLocal = "General"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormatLocal = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").NumberFormatLocal = "General"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormatLocal = "@"
$ws.Range("D22").Value = "472.82"
$ws.Range("D22").NumberFormatLocal = "General"
$ws.Range("E22").Value = "  +16.79%  "
$ws.Range("D23").NumberFormatLocal = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("D23").NumberFormatLocal = "General"
$ws.Range("E23").Value = "  +16.19%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormatLocal = "@"
$ws.Range("D25").Value = "86.51"
$ws.Range("D25").NumberFormatLocal = "General"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("D26").NumberFormatLocal = "@"
$ws.Range("D26").Value = "13.68"
$ws.Range("D26").NumberFormatLocal = "General"
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("D27").NumberFormatLocal = "@"
$ws.Range("D27").Value = "10.87"
$ws.Range("D27").NumberFormatLocal = "General"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +4.22%  "
$ws.Range("D29").NumberFormatLocal = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("D29").NumberFormatLocal = "General"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormatLocal = "@"
$ws.Range("D30").Value = "30.74"
$ws.Range("D30").NumberFormatLocal = "General"
$ws.Range("E30").Value = "  +5.04%  "
$ws.Range("D31").NumberFormatLocal = "@"
$ws.Range("D31").Value = "6.71"
$ws.Range("D31").NumberFormatLocal = "General"
$ws.Range("E31").Value = "  +3.41%  "
$ws.Range("D32").NumberFormatLocal = "@"
$ws.Range("D32").Value = "11.56"
$ws.Range("D32").NumberFormatLocal = "General"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormatLocal = "@"
$ws.Range("D33").Value = "580.23"
$ws.Range("D33").NumberFormatLocal = "General"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("D35").NumberFormatLocal = "@"
$ws.Range("D35").Value = "60.14"
$ws.Range("D35").NumberFormatLocal = "General"
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormatLocal = "@"
$ws.Range("D37").Value = "0.142"
$ws.Range("D37").NumberFormatLocal = "General"
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("D38").NumberFormatLocal = "@"
$ws.Range("D38").Value = "35.96"
$ws.Range("D38").NumberFormatLocal = "General"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormatLocal = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("D39").NumberFormatLocal = "General"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormatLocal = "@"
$ws.Range("D42").Value = "3.102.31"
$ws.Range("D42").NumberFormatLocal = "General"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").NumberFormatLocal = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").NumberFormatLocal = "General"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormatLocal = "@"
$ws.Range("D44").Value = "2.88"
$ws.Range("D44").NumberFormatLocal = "General"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").NumberFormatLocal = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").NumberFormatLocal = "General"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormatLocal = "@"
$ws.Range("D47").Value = "3.21"
$ws.Range("D47").NumberFormatLocal = "General"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormatLocal = "@"
$ws.Range("D48").Value = "0.134"
$ws.Range("D48").NumberFormatLocal = "General"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").NumberFormatLocal = "@"
$ws.Range("D50").Value = "137.63"
$ws.Range("D50").NumberFormatLocal = "General"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("D51").NumberFormatLocal = "@"
$ws.Range("D51").Value = "8.38"
$ws.Range("D51").NumberFormatLocal = "General"
$ws.Range("E51").Value = "  +3.05%  "
